$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: headers ---
# A1 loses its bold header style; becomes a plain identifier cell.
$ws.Range("A1").Value = "Sinvu#NA1"
$ws.Range("A1").ClearFormats()

# B1 already carries the bold/border header style; just replace its text.
$ws.Range("B1").Value = "Kills"

# C1:L1 need the same header style as B1 -> copy formats across.
$ws.Range("C1").Value = "Deaths"
$ws.Range("D1").Value = "Assists"
$ws.Range("E1").Value = "winRate"
$ws.Range("F1").Value = "avgCS"
$ws.Range("G1").Value = "GPM"
$ws.Range("H1").Value = "KP"
$ws.Range("I1").Value = "avgDPM"
$ws.Range("J1").Value = "avgAbusage"
$ws.Range("K1").Value = "visionScore"
$ws.Range("L1").Value = "numGames"

$ws.Range("B1").Copy()
$ws.Range("C1:L1").PasteSpecial(-4122)  # xlPasteFormats

# --- Row 2: data ---
# A2 gets the header style too (matches B1/C1 style), with the summoner name.
$ws.Range("A2").Value = "Yorick"
$ws.Range("B1").Copy()
$ws.Range("A2").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("B2").Value = 4
$ws.Range("C2").Value = 4
$ws.Range("D2").Value = 8
$ws.Range("E2").Value = 100
$ws.Range("F2").Value = 8
$ws.Range("G2").Value = 427

# H2 must stay a literal text string "0.0%", not get auto-converted to a percentage number.
$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "0.0%"
$ws.Range("H2").ClearFormats()

$ws.Range("I2").Value = 712
$ws.Range("J2").Value = 250
$ws.Range("K2").Value = 17
$ws.Range("L2").Value = 2

# --- Row 3 no longer exists in the new layout ---
$ws.Range("A3:B3").ClearContents()

$excel.CutCopyMode = $false
